$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The original column A (values 0,4,11,14 with bordered style) is a stray
# duplicate of the GENE column (old F). Delete that entire column; this
# shifts old columns B:F left to A:E, so the old B1 header
# ("QS_Astral_exact5") becomes the new A1 header automatically.
$ws.Range("A1").EntireColumn.Delete()
